$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.305.06"
$ws.Range("D3").Value = "1.610.13"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'213.31"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.88%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "'18.47"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.834.80"
$ws.Range("D13").Value = "1.600.59"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "'4.02"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "'0.515"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "26.273.52"
$ws.Range("D17").Value = "'62.23"
$ws.Range("E17").Value = "  +3.11%  "
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "'202.03"
$ws.Range("E20").Value = "  +0.47%  "
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "'9.33"
$ws.Range("E22").Value = "  +0.68%  "
$ws.Range("D23").Value = "'6.03"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").Value = "'1.89"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").Value = "'143.54"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'0.121"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "'15.25"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("E30").Value = "  +4.94%  "
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "'3.19"
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").Value = "1.162.13"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("D37").Value = "'0.0167"
$ws.Range("E37").Value = "  +2.17%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  +1.04%  "
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "1.746.16"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("D45").Value = "'92.45"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("E46").Value = "  +14.20%  "
$ws.Range("D47").Value = "'1.53"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("D49").Value = "'0.0508"
$ws.Range("E49").Value = "  +0.68%  "
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -0.09%  "
